# Apply the "Upload new version with timestamp" update to the Day Sale
# shortages report.
#
# Changes being applied (row 7 = ALKAPRESS PLUS 10/160MG 20 F.C. TABS.):
#   - الرصيد الحالي (current balance)  H7: "1:0"      -> "0:1"
#   - سعر البيع      (sale price)      P7: "102.0000" -> "153.0000"
#   - عدد التعاملات  (txn count)       Q7: "1:0"      -> "1:1"
#   - Row 13 total (sum of sell prices P7:P12) recalculated: 234.2 -> 285.2
#   - Row 14 "generated at" timestamp footer: 10:57 AM -> 10:59 AM

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 updates ---
$ws.Range("H7").Value = "0:1"
$ws.Range("Q7").Value = "1:1"

# P7 holds a numeric-looking string ("153.0000") that must stay stored as
# text (shared string), matching how the sheet already stores these
# formatted price strings. Writing it straight to .Value would make Excel
# coerce it into a real number, so instead we compute it via a formula and
# then collapse the formula down to a plain (text) value, which preserves
# both the original cell style and the text storage type.
$ws.Range("P7").Formula = '="153.0000"'
$ws.Range("P7").Copy()
$ws.Range("P7").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# --- Recalculate the displayed total for column P (sell price) ---
$ws.Range("P13").Value = 285.2

# --- Refresh the "generated at" timestamp footer ---
$ws.Range("A14").Value = "Wednesday, 17 September, 2025 10:59 AM"

$wb.Save()
